$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-11-19 Wednesday"; new = "2025-11-20 Thursday"},
    @{old = "504÷2="; new = "722÷7="},
    @{old = "133÷8="; new = "263÷5="},
    @{old = "570÷9="; new = "718÷9="},
    @{old = "474÷3="; new = "169÷2="},
    @{old = "433÷9="; new = "904÷2="},
    @{old = "572÷8="; new = "511÷8="},
    @{old = "459÷5="; new = "446÷3="},
    @{old = "367÷5="; new = "900÷4="},
    @{old = "455÷4="; new = "366÷6="},
    @{old = "292÷9="; new = "158÷7="},
    @{old = "435÷4="; new = "802÷6="},
    @{old = "286÷9="; new = "829÷5="},
    @{old = "996÷9="; new = "961÷5="},
    @{old = "668÷2="; new = "376÷4="},
    @{old = "677÷6="; new = "565÷3="},
    @{old = "631÷2="; new = "278÷9="},
    @{old = "909÷6="; new = "198÷6="},
    @{old = "999÷5="; new = "971÷4="},
    @{old = "255÷2="; new = "968÷5="},
    @{old = "177÷8="; new = "340÷6="},
    @{old = "218÷3="; new = "816÷3="},
    @{old = "800÷6="; new = "813÷3="},
    @{old = "222÷2="; new = "448÷6="},
    @{old = "170÷5="; new = "203÷6="},
    @{old = "165÷6="; new = "308÷3="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
